$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit removes one record (row 19) from the weekly price list; all
# subsequent rows shift up by one and the sheet's used range shrinks from
# A1:R57 to A1:R56.
$ws.Rows.Item(19).Delete()
